# B6-PowerPoint.pptx edit:
#  1. Re-style the three on-slide tables (slides 14-16) from the custom
#     "Table_0" style to the built-in "Medium Style 2 - Accent 1" style.
#  2. Swap the presentation's active theme colour scheme from the
#     "Integral" / "Red Violet" palette back to the stock "Office"
#     palette (this is what the slide master's theme - the one that is
#     actually painted behind every slide - ends up holding after the
#     edit).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table styles
# ---------------------------------------------------------------------
$newTableStyleId = "{53A7AA0B-39F4-4893-9410-620E9C6DE66A}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colours -> restore the stock "Office" palette
# ---------------------------------------------------------------------
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388    # dk2      44546A
    4  = 15132391   # lt2      E7E6E6
    5  = 13998939   # accent1  5B9BD5
    6  = 3243501    # accent2  ED7D31
    7  = 10855845   # accent3  A5A5A5
    8  = 49407      # accent4  FFC000
    9  = 12874308   # accent5  4472C4
    10 = 4697456    # accent6  70AD47
    11 = 12673797   # hlink    0563C1
    12 = 7491477    # folHlink 954F72
}

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i]
}
